# Refresh the "Price" (D) and "Volume(1h)" (E) columns of the cryptos sheet
# with the latest scraped figures, and splice in a newly-ranked coin
# (BabyDogeCoin) that bumped Aave/Cronos/Mantle down one row and pushed
# USDD off the bottom of the top-50 list.
#
# Note: several Price values look numeric (e.g. "65.00", "0.0844") but are
# authored as literal text in the source sheet (multi-dot big numbers like
# "26.510.14" make that obvious). A leading apostrophe forces Excel to store
# these as text instead of auto-converting/normalizing them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.510.14'
$ws.Range('E2').Value = '  -0.59%  '

$ws.Range('D3').Value = '1.627.30'
$ws.Range('E3').Value = '  -0.47%  '

$ws.Range('D5').Value = '''213.06'
$ws.Range('E5').Value = '  +0.07%  '

$ws.Range('D6').Value = '''0.502'
$ws.Range('E6').Value = '  +1.42%  '

$ws.Range('E7').Value = '  +0.12%  '

$ws.Range('E8').Value = '  -0.04%  '

$ws.Range('E9').Value = '  -1.59%  '

$ws.Range('E10').Value = '  -1.21%  '

$ws.Range('D11').Value = '''0.0844'
$ws.Range('E11').Value = '  +0.60%  '

$ws.Range('D12').Value = '1.852.95'

$ws.Range('D13').Value = '1.630.01'
$ws.Range('E13').Value = '  -0.40%  '

$ws.Range('E14').Value = '  +1.51%  '

$ws.Range('D15').Value = '''0.522'
$ws.Range('E15').Value = '  -0.63%  '

$ws.Range('D16').Value = '''65.00'
$ws.Range('E16').Value = '  +3.15%  '

$ws.Range('D17').Value = '26.542.13'
$ws.Range('E17').Value = '  -0.48%  '

$ws.Range('D18').Value = '0.0₃0740'
$ws.Range('E18').Value = '  -0.09%  '

$ws.Range('D19').Value = '''214.36'
$ws.Range('E19').Value = '  +2.82%  '

$ws.Range('E20').Value = '  +0.15%  '

$ws.Range('D21').Value = '''4.28'

$ws.Range('E22').Value = '  +1.52%  '

$ws.Range('D23').Value = '''9.28'
$ws.Range('E23').Value = '  -1.14%  '

$ws.Range('D24').Value = '''2.07'
$ws.Range('E24').Value = '  +8.69%  '

$ws.Range('D25').Value = '''148.32'
$ws.Range('E25').Value = '  +1.18%  '

$ws.Range('E27').Value = '  -0.16%  '

$ws.Range('E28').Value = '  +1.90%  '

$ws.Range('E29').Value = '  +0.87%  '

$ws.Range('D30').Value = '''0.0509'
$ws.Range('E30').Value = '  -2.19%  '

$ws.Range('E31').Value = '  -0.92%  '

$ws.Range('E32').Value = '  +3.15%  '

$ws.Range('E33').Value = '  -0.22%  '

$ws.Range('D34').Value = '1.238.11'
$ws.Range('E34').Value = '  +5.95%  '

$ws.Range('E35').Value = '  -0.02%  '

$ws.Range('E36').Value = '  -1.84%  '

$ws.Range('E37').Value = '  +4.04%  '

$ws.Range('E38').Value = '  +0.12%  '

$ws.Range('E39').Value = '  +0.68%  '

$ws.Range('E40').Value = '  -1.40%  '

$ws.Range('E41').Value = '  -1.87%  '

$ws.Range('D42').Value = '''0.798'
$ws.Range('E42').Value = '  +0.55%  '

$ws.Range('E43').Value = '  -0.77%  '

$ws.Range('D44').Value = '1.764.60'
$ws.Range('E44').Value = '  -0.67%  '

$ws.Range('D45').Value = '''93.01'
$ws.Range('E45').Value = '  +0.73%  '

$ws.Range('E46').Value = '  +2.51%  '

$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.0₆0104'
$ws.Range('E47').Value = '  +0.48%  '

$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '''54.87'
$ws.Range('E48').Value = '  +0.44%  '

$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.0509'
$ws.Range('E49').Value = '  -0.69%  '

$ws.Range('D50').Value = '''7.50'
$ws.Range('E50').Value = '  -0.38%  '

$ws.Range('B51').Value = 'Mantle'
$ws.Range('C51').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D51').Value = '''0.406'
$ws.Range('E51').Value = '  -0.80%  '
